# Update para U3 (ahora U2) y revisiones
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dynamic-array style formula for H72:H88, replacing the old
# "G*(1+0.963936)" shared formula. Each cell gets its own single-cell
# array formula (anchored to itself) that redistributes the
# (H - G) * C weighted average across rows 2:71 back onto G<row>.
for ($r = 72; $r -le 88; $r++) {
    $ws.Cells.Item($r, 8).FormulaArray = "=G$r+SUM((`$H`$2:`$H`$71 - `$G`$2:`$G`$71) * `$C`$2:`$C`$71) / SUM(`$C`$2:`$C`$71)"
}

# Column L (12) gets a custom width.
$ws.Columns.Item(12).ColumnWidth = 18.65

# Active cell moves to B6.
$ws.Range("B6").Select() | Out-Null
